$wb = $excel.ActiveWorkbook

$wsChange = $wb.Worksheets.Item("Änderungsgeschichte")
$wsRisk   = $wb.Worksheets.Item("Risiken")

# --- Änderungsgeschichte: fix wording in the v1.1 change description ---
$wsChange.Range("C5").Value = "Risiko 1 ist eingetreten: HW wird nicht rechtzeitig geliefert. Deshalb muss alternatives Testsetup evaluiert werden (Schaden: 25h).`nNeue Risiken hinzugefügt (4,6,7), Schätzungen angepasst, Beschreibungen erweitert"

# --- Änderungsgeschichte: append new row 7 for version 1.3 (review pass) ---
$loChange = $wsChange.ListObjects.Item(1)
$loChange.ListRows.Add() | Out-Null
$wsChange.Range("A7").Value = 40970
$wsChange.Range("B7").Value = "1.3"
$wsChange.Range("C7").Value = "Review, grammatikalische Korrekturen"
$wsChange.Range("D7").Value = "DT"

# --- Risiken: small grammar / typo fixes on existing risk texts ---
$wsRisk.Range("C4").Value = "Die Hardware für die Video Wall kann nicht  rechtzeitig geliefert werden."
$wsRisk.Range("G5").Value = "Scrum, früher Prototyp, kleine Demoprogramme zu Beginn des Projektes. Kontinuierliche Überprüfung der Projektplanung und eventuelle Anpassung. Verantwortung für aktuellen Projektplan an einem Teammitglied zuweisen."
$wsRisk.Range("C8").Value = "Die Auflösung der Video Wall ist für das Lesen der Bachelor Posters ungenügend."
$wsRisk.Range("H8").Value = "Darstellung der Poster in einer Grösse, die gut lesbar ist. Usability Tests müssen wiederholt werden. Alternative (2h Aufwand): Zu kleine Texte / Dokumente dürfen nicht hochgeladen / angezeigt werden (Constraint einführen)."

# --- Risiken: Risk 4 (Kinect Drehung) Sprint neu bumped from 1 to 2 ---
$wsRisk.Range("I7").Value = 2
